$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value. Values are written with a leading
# apostrophe to force text storage (preventing Excel from auto-converting
# numeric-looking strings like '0.998' into real numbers), then the cell
# style is reset to 'Normal' to avoid leaving a stray quote-prefix style.
$changes = [ordered]@{
    'D2' = '59.037.85'
    'E2' = '  -3.85%  '
    'D3' = '2.350.98'
    'E3' = '  -3.60%  '
    'D4' = '0.998'
    'E4' = '  +0.03%  '
    'D5' = '553.15'
    'E5' = '  -3.86%  '
    'D6' = '135.89'
    'E6' = '  -3.67%  '
    'E7' = '  +0.16%  '
    'D8' = '0.526'
    'E8' = '  -0.96%  '
    'D9' = '2.341.42'
    'E9' = '  -3.62%  '
    'E10' = '  -4.86%  '
    'E11' = '  -1.46%  '
    'E12' = '  -2.72%  '
    'D13' = '0.333'
    'E13' = '  -2.09%  '
    'D14' = '25.12'
    'E14' = '  -3.88%  '
    'D15' = '2.785.08'
    'E15' = '  -3.68%  '
    'D16' = '0.0000161'
    'E16' = '  -5.18%  '
    'D17' = '59.356.85'
    'E17' = '  -3.27%  '
    'B18' = 'WrappedEther'
    'C18' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D18' = '2.320.58'
    'E18' = '  -4.18%  '
    'B19' = 'Uniswap'
    'C19' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D19' = '8.02'
    'E19' = '  +10.62%  '
    'D20' = '10.31'
    'E20' = '  -2.75%  '
    'D21' = '318.37'
    'E21' = '  -1.95%  '
    'D22' = '3.97'
    'E22' = '  -2.08%  '
    'D23' = '5.95'
    'E23' = '  -1.81%  '
    'E24' = '  +0.04%  '
    'E25' = '  -7.70%  '
    'D26' = '63.56'
    'E26' = '  -2.22%  '
    'B27' = 'Aptos'
    'C27' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D27' = '7.97'
    'E27' = '  -11.51%  '
    'B28' = 'Bittensor'
    'C28' = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
    'D28' = '539.45'
    'E28' = '  -6.51%  '
    'D29' = '2.472.76'
    'E29' = '  -3.79%  '
    'D30' = '0.0₃0892'
    'E30' = '  -2.26%  '
    'D31' = '7.79'
    'E31' = '  -1.09%  '
    'D32' = '1.27'
    'E32' = '  -5.71%  '
    'D33' = '1.74'
    'E33' = '  -6.11%  '
    'D34' = '0.127'
    'E34' = '  -4.28%  '
    'E35' = '  -0.27%  '
    'D36' = '150.52'
    'E36' = '  -0.88%  '
    'D37' = '1.38'
    'E37' = '  -0.15%  '
    'D38' = '0.361'
    'E38' = '  -2.42%  '
    'D39' = '4.46'
    'E39' = '  -3.92%  '
    'D40' = '17.92'
    'E40' = '  -2.14%  '
    'D41' = '4.94'
    'E41' = '  -3.68%  '
    'E42' = '  +0.00%  '
    'D43' = '41.04'
    'E43' = '  -1.64%  '
    'D44' = '1.61'
    'E44' = '  -3.60%  '
    'D45' = '2.29'
    'E45' = '  -2.44%  '
    'D46' = '0.0₆0279'
    'E46' = '  -5.06%  '
    'D47' = '137.12'
    'E47' = '  -3.60%  '
    'D48' = '3.46'
    'E48' = '  -2.42%  '
    'D49' = '0.576'
    'E49' = '  -3.22%  '
    'D50' = '0.0492'
    'E50' = '  -3.05%  '
    'D51' = '18.69'
    'E51' = '  -4.40%  '
}

foreach ($cellRef in $changes.Keys) {
    $ws.Range($cellRef).Value = "'" + $changes[$cellRef]
    $ws.Range($cellRef).Style = "Normal"
}
